# "Finally got tests working"
#
# The fixture's Sheet1 header row originally was:
#   A=Last Name, B=First Name, C=Teen Name, D=Grade, E=Gender, F=School,
#   G=Teen Email, H=Parent Email Address(es), I=Cardinal Gibbons HS Group
#
# The three columns that are not needed by the spec (Teen Name, Gender,
# School) are removed entirely, collapsing the sheet down to:
#   A=Last Name, B=First Name, C=Grade, D=Teen Email,
#   E=Parent Email Address(es), F=Cardinal Gibbons HS Group
#
# and the sheet's selection is moved off the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete the unwanted columns right-to-left so earlier deletes don't
# invalidate the column letters used by later ones.
$ws.Range("F1").EntireColumn.Delete()   # School
$ws.Range("E1").EntireColumn.Delete()   # Gender
$ws.Range("C1").EntireColumn.Delete()   # Teen Name

# Match the saved selection/active cell from the edited workbook.
$ws.Range("E31").Select()
